# Add a new "docente" (teacher) record as row 3 on the "docentes" sheet.
# Columns: A Nombre, B Apellido, C Documento, D Usuario, E Correo,
#          F Telefono, G Contraseña, H Perfil, I Departamento, J Ciudad, K Universidad

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("docentes")

$row = 3

# Plain text fields - safe to assign directly, Excel stores them as text.
$ws.Cells.Item($row, 1).Value2 = "don"
$ws.Cells.Item($row, 2).Value2 = "papaspapas"

# "Documento" looks purely numeric - force it to be stored as text (not a
# number) the same way the rest of the sheet does, without leaving a
# leftover "quote prefix" cell style behind: build it with a text formula
# and then convert the formula to a plain value in place.
$ws.Cells.Item($row, 3).Formula = '="1032019684"'
$ws.Cells.Item($row, 3).Copy()
$ws.Cells.Item($row, 3).PasteSpecial(-4163)

$ws.Cells.Item($row, 4).Value2 = "donpapas35"
$ws.Cells.Item($row, 5).Value2 = "donpapas35@salchipapas.com"

# "Telefono" is also purely numeric text - same trick as "Documento".
$ws.Cells.Item($row, 6).Formula = '="3535353535"'
$ws.Cells.Item($row, 6).Copy()
$ws.Cells.Item($row, 6).PasteSpecial(-4163)

$ws.Cells.Item($row, 7).Value2 = "donpapa"
$ws.Cells.Item($row, 8).Value2 = "Docente"
$ws.Cells.Item($row, 9).Value2 = "Antioquia"
$ws.Cells.Item($row, 10).Value2 = "Medellín"
$ws.Cells.Item($row, 11).Value2 = "El poli"
